# The workbook is a Selenium/Java test-data spreadsheet (SauceDemo login
# creds). The author re-saved the file on their own machine after filling
# in two previously-blank "extra login" rows on the "Login Details" sheet
# with their own surname/name ("Cele" / "Nkosi"), then scrolled/zoomed the
# view before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Details")

# --- Data edits -----------------------------------------------------
# Row 5, column B previously held a stray "secret_sauce" value; replace it
# with the surname "Cele".
$ws.Range("B5").Value = "Cele"

# Row 7, column B previously held a stray "secret_sauce" value; replace it
# with the name "Nkosi".
$ws.Range("B7").Value = "Nkosi"

# --- View state (best-effort, matches the author's last saved view) -----
$ws.Activate()
$ws.Range("E612").Select() | Out-Null
$excel.ActiveWindow.Zoom = 180
